$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.353747141763602
$ws.Range("D2").Value = 0.1895623177382468

$ws.Range("C3").Value = 0.4499373804056233
$ws.Range("D3").Value = 0.6571577243561531

$ws.Range("C4").Value = 1.713831849076387
$ws.Range("D4").Value = 0.1006167898480195

$ws.Range("C5").Value = 0.04204420896542845
$ws.Range("D5").Value = 0.9668427098033123

$ws.Range("C6").Value = -0.858360284277333
$ws.Range("D6").Value = 0.3999491017450074

$ws.Range("C7").Value = 0.5718977486511369
$ws.Range("D7").Value = 0.5731889804864032

$ws.Range("C8").Value = -1.278206925987671
$ws.Range("D8").Value = 0.2144982848023247

$ws.Range("C9").Value = 1.265188250742409
$ws.Range("D9").Value = 0.2190409061330811

$ws.Range("C10").Value = -0.583173122523622
$ws.Range("D10").Value = 0.565707573782456

$ws.Range("C11").Value = -1.550516962648961
$ws.Range("D11").Value = 0.1352847225270906
